$wb = $excel.ActiveWorkbook

# --- Sheet1 (BD_Times): add rows 356-397 ---
$ws1 = $wb.Worksheets.Item("BD_Times")
$arr1 = New-Object 'object[,]' 42,9
$arr1[0,0] = 'Botafogo'
$arr1[0,1] = 1
$arr1[0,2] = 1
$arr1[0,3] = 1
$arr1[0,4] = 1
$arr1[0,5] = 3
$arr1[0,6] = 1
$arr1[0,7] = 7
$arr1[0,8] = 3
$arr1[1,0] = 'Internacional'
$arr1[1,1] = 0
$arr1[1,2] = 1
$arr1[1,3] = 1
$arr1[1,4] = 1
$arr1[1,5] = 1
$arr1[1,6] = 3
$arr1[1,7] = 3
$arr1[1,8] = 7
$arr1[2,0] = 'Atletico MG'
$arr1[2,1] = 1
$arr1[2,2] = 1
$arr1[2,3] = 0
$arr1[2,4] = 0
$arr1[2,5] = 1
$arr1[2,6] = 0
$arr1[2,7] = 2
$arr1[2,8] = 9
$arr1[3,0] = 'Bahia'
$arr1[3,1] = 0
$arr1[3,2] = 0
$arr1[3,3] = 1
$arr1[3,4] = 0
$arr1[3,5] = 0
$arr1[3,6] = 1
$arr1[3,7] = 9
$arr1[3,8] = 2
$arr1[4,0] = 'America'
$arr1[4,1] = 1
$arr1[4,2] = 0
$arr1[4,3] = 1
$arr1[4,4] = 0
$arr1[4,5] = 0
$arr1[4,6] = 1
$arr1[4,7] = 7
$arr1[4,8] = 3
$arr1[5,0] = 'Goias'
$arr1[5,1] = 0
$arr1[5,2] = 1
$arr1[5,3] = 0
$arr1[5,4] = 0
$arr1[5,5] = 1
$arr1[5,6] = 0
$arr1[5,7] = 3
$arr1[5,8] = 7
$arr1[6,0] = 'Corinthians'
$arr1[6,1] = 1
$arr1[6,2] = 1
$arr1[6,3] = 1
$arr1[6,4] = 1
$arr1[6,5] = 3
$arr1[6,6] = 1
$arr1[6,7] = 4
$arr1[6,8] = 2
$arr1[7,0] = 'Coritiba'
$arr1[7,1] = 0
$arr1[7,2] = 1
$arr1[7,3] = 1
$arr1[7,4] = 1
$arr1[7,5] = 1
$arr1[7,6] = 3
$arr1[7,7] = 2
$arr1[7,8] = 4
$arr1[8,0] = 'Gremio'
$arr1[8,1] = 1
$arr1[8,2] = 1
$arr1[8,3] = 1
$arr1[8,4] = 1
$arr1[8,5] = 2
$arr1[8,6] = 1
$arr1[8,7] = 3
$arr1[8,8] = 8
$arr1[9,0] = 'Fluminense'
$arr1[9,1] = 0
$arr1[9,2] = 1
$arr1[9,3] = 1
$arr1[9,4] = 1
$arr1[9,5] = 1
$arr1[9,6] = 2
$arr1[9,7] = 8
$arr1[9,8] = 3
$arr1[10,0] = 'Flamengo'
$arr1[10,1] = 1
$arr1[10,2] = 1
$arr1[10,3] = 1
$arr1[10,4] = 1
$arr1[10,5] = 1
$arr1[10,6] = 1
$arr1[10,7] = 8
$arr1[10,8] = 2
$arr1[11,0] = 'Sao Paulo'
$arr1[11,1] = 0
$arr1[11,2] = 1
$arr1[11,3] = 1
$arr1[11,4] = 1
$arr1[11,5] = 1
$arr1[11,6] = 1
$arr1[11,7] = 2
$arr1[11,8] = 8
$arr1[12,0] = 'Fortaleza'
$arr1[12,1] = 1
$arr1[12,2] = 1
$arr1[12,3] = 0
$arr1[12,4] = 0
$arr1[12,5] = 4
$arr1[12,6] = 0
$arr1[12,7] = 10
$arr1[12,8] = 6
$arr1[13,0] = 'Santos'
$arr1[13,1] = 0
$arr1[13,2] = 0
$arr1[13,3] = 1
$arr1[13,4] = 0
$arr1[13,5] = 0
$arr1[13,6] = 4
$arr1[13,7] = 6
$arr1[13,8] = 10
$arr1[14,0] = 'Palmeiras'
$arr1[14,1] = 1
$arr1[14,2] = 1
$arr1[14,3] = 0
$arr1[14,4] = 0
$arr1[14,5] = 1
$arr1[14,6] = 0
$arr1[14,7] = 3
$arr1[14,8] = 5
$arr1[15,0] = 'Cruzeiro'
$arr1[15,1] = 0
$arr1[15,2] = 0
$arr1[15,3] = 1
$arr1[15,4] = 0
$arr1[15,5] = 0
$arr1[15,6] = 1
$arr1[15,7] = 5
$arr1[15,8] = 3
$arr1[16,0] = 'Bragantino'
$arr1[16,1] = 1
$arr1[16,2] = 1
$arr1[16,3] = 1
$arr1[16,4] = 1
$arr1[16,5] = 1
$arr1[16,6] = 1
$arr1[16,7] = 7
$arr1[16,8] = 6
$arr1[17,0] = 'Vasco'
$arr1[17,1] = 0
$arr1[17,2] = 1
$arr1[17,3] = 1
$arr1[17,4] = 1
$arr1[17,5] = 1
$arr1[17,6] = 1
$arr1[17,7] = 6
$arr1[17,8] = 7
$arr1[18,0] = 'Athletico PR'
$arr1[18,1] = 1
$arr1[18,2] = 1
$arr1[18,3] = 0
$arr1[18,4] = 0
$arr1[18,5] = 2
$arr1[18,6] = 0
$arr1[18,7] = 5
$arr1[18,8] = 3
$arr1[19,0] = 'Cuiaba'
$arr1[19,1] = 0
$arr1[19,2] = 0
$arr1[19,3] = 1
$arr1[19,4] = 0
$arr1[19,5] = 0
$arr1[19,6] = 2
$arr1[19,7] = 3
$arr1[19,8] = 5
$arr1[20,0] = 'Internacional'
$arr1[20,1] = 1
$arr1[20,2] = 0
$arr1[20,3] = 1
$arr1[20,4] = 0
$arr1[20,5] = 0
$arr1[20,6] = 1
$arr1[20,7] = 10
$arr1[20,8] = 6
$arr1[21,0] = 'Fortaleza'
$arr1[21,1] = 0
$arr1[21,2] = 1
$arr1[21,3] = 0
$arr1[21,4] = 0
$arr1[21,5] = 1
$arr1[21,6] = 0
$arr1[21,7] = 6
$arr1[21,8] = 10
$arr1[22,0] = 'Sao Paulo'
$arr1[22,1] = 1
$arr1[22,2] = 0
$arr1[22,3] = 0
$arr1[22,4] = 0
$arr1[22,5] = 0
$arr1[22,6] = 0
$arr1[22,7] = 5
$arr1[22,8] = 2
$arr1[23,0] = 'Botafogo'
$arr1[23,1] = 0
$arr1[23,2] = 0
$arr1[23,3] = 0
$arr1[23,4] = 0
$arr1[23,5] = 0
$arr1[23,6] = 0
$arr1[23,7] = 2
$arr1[23,8] = 5
$arr1[24,0] = 'Cuiaba'
$arr1[24,1] = 1
$arr1[24,2] = 0
$arr1[24,3] = 1
$arr1[24,4] = 0
$arr1[24,5] = 0
$arr1[24,6] = 2
$arr1[24,7] = 9
$arr1[24,8] = 2
$arr1[25,0] = 'Palmeiras'
$arr1[25,1] = 0
$arr1[25,2] = 1
$arr1[25,3] = 0
$arr1[25,4] = 0
$arr1[25,5] = 2
$arr1[25,6] = 0
$arr1[25,7] = 2
$arr1[25,8] = 9
$arr1[26,0] = 'Fluminense'
$arr1[26,1] = 1
$arr1[26,2] = 1
$arr1[26,3] = 1
$arr1[26,4] = 1
$arr1[26,5] = 3
$arr1[26,6] = 1
$arr1[26,7] = 7
$arr1[26,8] = 1
$arr1[27,0] = 'America'
$arr1[27,1] = 0
$arr1[27,2] = 1
$arr1[27,3] = 1
$arr1[27,4] = 1
$arr1[27,5] = 1
$arr1[27,6] = 3
$arr1[27,7] = 1
$arr1[27,8] = 7
$arr1[28,0] = 'Cruzeiro'
$arr1[28,1] = 1
$arr1[28,2] = 1
$arr1[28,3] = 1
$arr1[28,4] = 1
$arr1[28,5] = 1
$arr1[28,6] = 1
$arr1[28,7] = 10
$arr1[28,8] = 3
$arr1[29,0] = 'Corinthians'
$arr1[29,1] = 0
$arr1[29,2] = 1
$arr1[29,3] = 1
$arr1[29,4] = 1
$arr1[29,5] = 1
$arr1[29,6] = 1
$arr1[29,7] = 3
$arr1[29,8] = 10
$arr1[30,0] = 'Vasco'
$arr1[30,1] = 1
$arr1[30,2] = 1
$arr1[30,3] = 0
$arr1[30,4] = 0
$arr1[30,5] = 1
$arr1[30,6] = 0
$arr1[30,7] = 3
$arr1[30,8] = 7
$arr1[31,0] = 'Atletico MG'
$arr1[31,1] = 0
$arr1[31,2] = 0
$arr1[31,3] = 1
$arr1[31,4] = 0
$arr1[31,5] = 0
$arr1[31,6] = 1
$arr1[31,7] = 7
$arr1[31,8] = 3
$arr1[32,0] = 'Bahia'
$arr1[32,1] = 1
$arr1[32,2] = 1
$arr1[32,3] = 0
$arr1[32,4] = 0
$arr1[32,5] = 4
$arr1[32,6] = 0
$arr1[32,7] = 4
$arr1[32,8] = 7
$arr1[33,0] = 'Bragantino'
$arr1[33,1] = 0
$arr1[33,2] = 0
$arr1[33,3] = 1
$arr1[33,4] = 0
$arr1[33,5] = 0
$arr1[33,6] = 4
$arr1[33,7] = 7
$arr1[33,8] = 4
$arr1[34,0] = 'Coritiba'
$arr1[34,1] = 1
$arr1[34,2] = 1
$arr1[34,3] = 1
$arr1[34,4] = 1
$arr1[34,5] = 2
$arr1[34,6] = 3
$arr1[34,7] = 4
$arr1[34,8] = 5
$arr1[35,0] = 'Flamengo'
$arr1[35,1] = 0
$arr1[35,2] = 1
$arr1[35,3] = 1
$arr1[35,4] = 1
$arr1[35,5] = 3
$arr1[35,6] = 2
$arr1[35,7] = 5
$arr1[35,8] = 4
$arr1[36,0] = 'Santos'
$arr1[36,1] = 1
$arr1[36,2] = 1
$arr1[36,3] = 1
$arr1[36,4] = 1
$arr1[36,5] = 2
$arr1[36,6] = 1
$arr1[36,7] = 9
$arr1[36,8] = 4
$arr1[37,0] = 'Gremio'
$arr1[37,1] = 0
$arr1[37,2] = 1
$arr1[37,3] = 1
$arr1[37,4] = 1
$arr1[37,5] = 1
$arr1[37,6] = 2
$arr1[37,7] = 4
$arr1[37,8] = 9
$arr1[38,0] = 'Goias'
$arr1[38,1] = 1
$arr1[38,2] = 1
$arr1[38,3] = 1
$arr1[38,4] = 1
$arr1[38,5] = 1
$arr1[38,6] = 1
$arr1[38,7] = 4
$arr1[38,8] = 7
$arr1[39,0] = 'Athletico PR'
$arr1[39,1] = 0
$arr1[39,2] = 1
$arr1[39,3] = 1
$arr1[39,4] = 1
$arr1[39,5] = 1
$arr1[39,6] = 1
$arr1[39,7] = 7
$arr1[39,8] = 4
$arr1[40,0] = 'Flamengo'
$arr1[40,1] = 1
$arr1[40,2] = 0
$arr1[40,3] = 0
$arr1[40,4] = 0
$arr1[40,5] = 0
$arr1[40,6] = 0
$arr1[40,7] = 6
$arr1[40,8] = 3
$arr1[41,0] = 'Internacional'
$arr1[41,1] = 0
$arr1[41,2] = 0
$arr1[41,3] = 0
$arr1[41,4] = 0
$arr1[41,5] = 0
$arr1[41,6] = 0
$arr1[41,7] = 3
$arr1[41,8] = 6
$ws1.Range("A356:I397").Value = $arr1

# --- Sheet2 (BD_Jogo): add rows 179-199 ---
$ws2 = $wb.Worksheets.Item("BD_Jogo")
$arr2 = New-Object 'object[,]' 21,5
$arr2[0,0] = 1
$arr2[0,1] = 4
$arr2[0,2] = 10
$arr2[0,3] = 'Botafogo'
$arr2[0,4] = 'Internacional'
$arr2[1,0] = 0
$arr2[1,1] = 1
$arr2[1,2] = 11
$arr2[1,3] = 'Atletico MG'
$arr2[1,4] = 'Bahia'
$arr2[2,0] = 0
$arr2[2,1] = 1
$arr2[2,2] = 10
$arr2[2,3] = 'America'
$arr2[2,4] = 'Goias'
$arr2[3,0] = 1
$arr2[3,1] = 4
$arr2[3,2] = 6
$arr2[3,3] = 'Corinthians'
$arr2[3,4] = 'Coritiba'
$arr2[4,0] = 1
$arr2[4,1] = 3
$arr2[4,2] = 11
$arr2[4,3] = 'Gremio'
$arr2[4,4] = 'Fluminense'
$arr2[5,0] = 1
$arr2[5,1] = 2
$arr2[5,2] = 10
$arr2[5,3] = 'Flamengo'
$arr2[5,4] = 'Sao Paulo'
$arr2[6,0] = 0
$arr2[6,1] = 4
$arr2[6,2] = 16
$arr2[6,3] = 'Fortaleza'
$arr2[6,4] = 'Santos'
$arr2[7,0] = 0
$arr2[7,1] = 1
$arr2[7,2] = 8
$arr2[7,3] = 'Palmeiras'
$arr2[7,4] = 'Cruzeiro'
$arr2[8,0] = 1
$arr2[8,1] = 2
$arr2[8,2] = 13
$arr2[8,3] = 'Bragantino'
$arr2[8,4] = 'Vasco'
$arr2[9,0] = 0
$arr2[9,1] = 2
$arr2[9,2] = 8
$arr2[9,3] = 'Athletico PR'
$arr2[9,4] = 'Cuiaba'
$arr2[10,0] = 0
$arr2[10,1] = 1
$arr2[10,2] = 16
$arr2[10,3] = 'Internacional'
$arr2[10,4] = 'Fortaleza'
$arr2[11,0] = 0
$arr2[11,1] = 0
$arr2[11,2] = 7
$arr2[11,3] = 'Sao Paulo'
$arr2[11,4] = 'Botafogo'
$arr2[12,0] = 0
$arr2[12,1] = 2
$arr2[12,2] = 11
$arr2[12,3] = 'Cuiaba'
$arr2[12,4] = 'Palmeiras'
$arr2[13,0] = 1
$arr2[13,1] = 4
$arr2[13,2] = 8
$arr2[13,3] = 'Fluminense'
$arr2[13,4] = 'America'
$arr2[14,0] = 1
$arr2[14,1] = 2
$arr2[14,2] = 13
$arr2[14,3] = 'Cruzeiro'
$arr2[14,4] = 'Corinthians'
$arr2[15,0] = 0
$arr2[15,1] = 1
$arr2[15,2] = 10
$arr2[15,3] = 'Vasco'
$arr2[15,4] = 'Atletico MG'
$arr2[16,0] = 0
$arr2[16,1] = 4
$arr2[16,2] = 11
$arr2[16,3] = 'Bahia'
$arr2[16,4] = 'Bragantino'
$arr2[17,0] = 1
$arr2[17,1] = 5
$arr2[17,2] = 9
$arr2[17,3] = 'Coritiba'
$arr2[17,4] = 'Flamengo'
$arr2[18,0] = 1
$arr2[18,1] = 3
$arr2[18,2] = 13
$arr2[18,3] = 'Santos'
$arr2[18,4] = 'Gremio'
$arr2[19,0] = 1
$arr2[19,1] = 2
$arr2[19,2] = 11
$arr2[19,3] = 'Goias'
$arr2[19,4] = 'Athletico PR'
$arr2[20,0] = 0
$arr2[20,1] = 0
$arr2[20,2] = 9
$arr2[20,3] = 'Flamengo'
$arr2[20,4] = 'Internacional'
$ws2.Range("A179:E199").Value = $arr2
